# Refresh scraped event data (per the commit "Update gh-pages to output
# generated at 456a3b4"):
#  - "want to go" counts (column F) bump for several already-listed events
#  - a newly discovered event ("合肥·Holic动漫游戏展") is inserted as row 13 on
#    both the "展览" (Exhibitions) and "全部类型" (All types) sheets, pushing all
#    later rows down by one; a few of those shifted rows also get small bumps
#    to their "want to go" counts
$wb = $excel.ActiveWorkbook

# ---------------- Sheet: 展览 ----------------
$ws1 = $wb.Worksheets.Item('展览')

# Updated "want to go" counts for existing (unmoved) rows
$ws1.Range('F2').Value = 724
$ws1.Range('F3').Value = 582
$ws1.Range('F4').Value = 551
$ws1.Range('F7').Value = 62
$ws1.Range('F9').Value = 53
$ws1.Range('F11').Value = 4730
$ws1.Range('F12').Value = 4502

# Insert a new row so rows 13..16 shift down to 14..17
$ws1.Rows.Item(13).Insert()

# The blank inserted row loses the bordered/bold "index column" style that
# every other row in column A carries; copy it back from the row below so the
# new row 13 matches the rest of column A.
$ws1.Range('A14').Copy()
$ws1.Range('A13').PasteSpecial(-4122)

# Write final contents for the new row and every shifted row
$ws1.Range('A13').Value = 12
$ws1.Range('B13').Value = '''2024-10-04'
$ws1.Range('C13').Value = '''合肥·Holic动漫游戏展'
$ws1.Range('D13').Value = '''庐州大道800号 合肥融创茂'
$ws1.Range('E13').Value = '''2024.10.04 10:00-10.06 17:00'
$ws1.Range('F13').Value = 0
$ws1.Range('G13').Value = 55
$ws1.Range('H13').Value = '''https://show.bilibili.com/platform/detail.html?id=92061'
$ws1.Range('I13').Value = '''//i1.hdslb.com/bfs/openplatform/202409/AZ0LsUce1725522015668.jpeg'

$ws1.Range('A14').Value = 13
$ws1.Range('B14').Value = '''2024-10-04'
$ws1.Range('C14').Value = '''合肥·乐帮•崩原铁绝only同人首展'
$ws1.Range('D14').Value = '''丹霞路488号金星商业城三楼 迷鹿轰趴'
$ws1.Range('E14').Value = '''2024.10.04 10:00-10.05 16:30'
$ws1.Range('F14').Value = 19
$ws1.Range('G14').Value = 58
$ws1.Range('H14').Value = '''https://show.bilibili.com/platform/detail.html?id=91524'
$ws1.Range('I14').Value = '''//i2.hdslb.com/bfs/openplatform/202408/739I7YRr1724912450704.png'

$ws1.Range('A15').Value = 14
$ws1.Range('B15').Value = '''2024-10-06'
$ws1.Range('C15').Value = '''合肥·星月动漫游戏展'
$ws1.Range('D15').Value = '''灵石路与皇藏峪路交叉口西南10米安徽百事兴电气有限公司院内2栋厂房2层 兄弟篮球俱乐部'
$ws1.Range('E15').Value = '''2024.10.06 10:00-10.06 17:00'
$ws1.Range('F15').Value = 1
$ws1.Range('G15').Value = 45
$ws1.Range('H15').Value = '''https://show.bilibili.com/platform/detail.html?id=91958'
$ws1.Range('I15').Value = '''//i2.hdslb.com/bfs/openplatform/202409/mgB8U6bN1725361649767.jpeg'

$ws1.Range('A16').Value = 15
$ws1.Range('B16').Value = '''2024-10-06'
$ws1.Range('C16').Value = '''合肥·首届火影忍者同人only'
$ws1.Range('D16').Value = '''长江东路金太阳家具广场南门二楼 优极篮球馆'
$ws1.Range('E16').Value = '''2024.10.06 09:30-10.06 17:30'
$ws1.Range('F16').Value = 26
$ws1.Range('G16').Value = 75
$ws1.Range('H16').Value = '''https://show.bilibili.com/platform/detail.html?id=91658'
$ws1.Range('I16').Value = '''//i0.hdslb.com/bfs/openplatform/202408/f8ylbskH1725027552569.jpeg'

$ws1.Range('A17').Value = 16
$ws1.Range('B17').Value = '''2024-10-26'
$ws1.Range('C17').Value = '''合肥·W·A第五人格同人only2.0'
$ws1.Range('D17').Value = '''莲花路与石门路交口西北角（尚泽大都会B座四楼） 格律诗婚礼艺术中心(经开店)'
$ws1.Range('E17').Value = '''2024.10.26 09:30-10.26 17:00'
$ws1.Range('F17').Value = 160
$ws1.Range('G17').Value = 68
$ws1.Range('H17').Value = '''https://show.bilibili.com/platform/detail.html?id=91123'
$ws1.Range('I17').Value = '''//i2.hdslb.com/bfs/openplatform/202408/YqXHTFM81724066565119.png'

# ---------------- Sheet: 全部类型 ----------------
$ws4 = $wb.Worksheets.Item('全部类型')

# Updated "want to go" counts for existing (unmoved) rows
$ws4.Range('F2').Value = 724
$ws4.Range('F3').Value = 582
$ws4.Range('F4').Value = 551
$ws4.Range('F7').Value = 62
$ws4.Range('F9').Value = 53
$ws4.Range('F11').Value = 4730
$ws4.Range('F12').Value = 4502

# Insert a new row so rows 13..19 shift down to 14..20
$ws4.Rows.Item(13).Insert()

# The blank inserted row loses the bordered/bold "index column" style that
# every other row in column A carries; copy it back from the row below so the
# new row 13 matches the rest of column A.
$ws4.Range('A14').Copy()
$ws4.Range('A13').PasteSpecial(-4122)

# Write final contents for the new row and every shifted row
$ws4.Range('A13').Value = 12
$ws4.Range('B13').Value = '''2024-10-04'
$ws4.Range('C13').Value = '''合肥·Holic动漫游戏展'
$ws4.Range('D13').Value = '''庐州大道800号 合肥融创茂'
$ws4.Range('E13').Value = '''2024.10.04 10:00-10.06 17:00'
$ws4.Range('F13').Value = 0
$ws4.Range('G13').Value = 55
$ws4.Range('H13').Value = '''https://show.bilibili.com/platform/detail.html?id=92061'
$ws4.Range('I13').Value = '''//i1.hdslb.com/bfs/openplatform/202409/AZ0LsUce1725522015668.jpeg'

$ws4.Range('A14').Value = 13
$ws4.Range('B14').Value = '''2024-10-04'
$ws4.Range('C14').Value = '''合肥·乐帮•崩原铁绝only同人首展'
$ws4.Range('D14').Value = '''丹霞路488号金星商业城三楼 迷鹿轰趴'
$ws4.Range('E14').Value = '''2024.10.04 10:00-10.05 16:30'
$ws4.Range('F14').Value = 19
$ws4.Range('G14').Value = 58
$ws4.Range('H14').Value = '''https://show.bilibili.com/platform/detail.html?id=91524'
$ws4.Range('I14').Value = '''//i2.hdslb.com/bfs/openplatform/202408/739I7YRr1724912450704.png'

$ws4.Range('A15').Value = 14
$ws4.Range('B15').Value = '''2024-10-06'
$ws4.Range('C15').Value = '''合肥·星月动漫游戏展'
$ws4.Range('D15').Value = '''灵石路与皇藏峪路交叉口西南10米安徽百事兴电气有限公司院内2栋厂房2层 兄弟篮球俱乐部'
$ws4.Range('E15').Value = '''2024.10.06 10:00-10.06 17:00'
$ws4.Range('F15').Value = 1
$ws4.Range('G15').Value = 45
$ws4.Range('H15').Value = '''https://show.bilibili.com/platform/detail.html?id=91958'
$ws4.Range('I15').Value = '''//i2.hdslb.com/bfs/openplatform/202409/mgB8U6bN1725361649767.jpeg'

$ws4.Range('A16').Value = 15
$ws4.Range('B16').Value = '''2024-10-06'
$ws4.Range('C16').Value = '''合肥·首届火影忍者同人only'
$ws4.Range('D16').Value = '''长江东路金太阳家具广场南门二楼 优极篮球馆'
$ws4.Range('E16').Value = '''2024.10.06 09:30-10.06 17:30'
$ws4.Range('F16').Value = 26
$ws4.Range('G16').Value = 75
$ws4.Range('H16').Value = '''https://show.bilibili.com/platform/detail.html?id=91658'
$ws4.Range('I16').Value = '''//i0.hdslb.com/bfs/openplatform/202408/f8ylbskH1725027552569.jpeg'

$ws4.Range('A17').Value = 16
$ws4.Range('B17').Value = '''2024-10-26'
$ws4.Range('C17').Value = '''合肥·W·A第五人格同人only2.0'
$ws4.Range('D17').Value = '''莲花路与石门路交口西北角（尚泽大都会B座四楼） 格律诗婚礼艺术中心(经开店)'
$ws4.Range('E17').Value = '''2024.10.26 09:30-10.26 17:00'
$ws4.Range('F17').Value = 160
$ws4.Range('G17').Value = 68
$ws4.Range('H17').Value = '''https://show.bilibili.com/platform/detail.html?id=91123'
$ws4.Range('I17').Value = '''//i2.hdslb.com/bfs/openplatform/202408/YqXHTFM81724066565119.png'

$ws4.Range('A18').Value = 17
$ws4.Range('B18').Value = '''2024-10-26'
$ws4.Range('C18').Value = '''合肥·《四月是你的谎言》—“公生”与“薰”的钢琴小提琴唯美经典音乐集'
$ws4.Range('D18').Value = '''徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院'
$ws4.Range('E18').Value = '''2024.10.26 19:30-10.26 21:00'
$ws4.Range('F18').Value = 66
$ws4.Range('G18').Value = 80
$ws4.Range('H18').Value = '''https://show.bilibili.com/platform/detail.html?id=90322'
$ws4.Range('I18').Value = '''//i2.hdslb.com/bfs/openplatform/202408/BiVgXUKH1722824304648.jpeg'

$ws4.Range('A19').Value = 18
$ws4.Range('B19').Value = '''2024-11-09'
$ws4.Range('C19').Value = '''合肥·一生必听的钢琴曲—“从巴赫 · 莫扎特到肖邦 · 李斯特”钢琴圣手谭小棠独奏音乐会'
$ws4.Range('D19').Value = '''徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院'
$ws4.Range('E19').Value = '''2024.11.09 19:30-11.09 21:00'
$ws4.Range('F19').Value = 5
$ws4.Range('G19').Value = 64
$ws4.Range('H19').Value = '''https://show.bilibili.com/platform/detail.html?id=90593'
$ws4.Range('I19').Value = '''//i2.hdslb.com/bfs/openplatform/202408/SYfLxnO21723442234232.jpeg'

$ws4.Range('A20').Value = 19
$ws4.Range('B20').Value = '''2024-12-07'
$ws4.Range('C20').Value = '''合肥·一生必听的古典系列《钟》—超技钢琴曲炫彩音乐会'
$ws4.Range('D20').Value = '''徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院'
$ws4.Range('E20').Value = '''2024.12.07 19:30-12.07 21:00'
$ws4.Range('F20').Value = 0
$ws4.Range('G20').Value = 56
$ws4.Range('H20').Value = '''https://show.bilibili.com/platform/detail.html?id=91608'
$ws4.Range('I20').Value = '''//i0.hdslb.com/bfs/openplatform/202408/wiLiWoeM1725005636569.jpeg'

